# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data block (row 277),
# shifting all existing records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 277; this shifts rows 277:357 down to 278:358
# and carries the row's number formatting (e.g. the date style on column D).
$ws.Rows("277:277").Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A277").Value = 10
$ws.Range("B277").Value = "Vega Modelo de Temuco"
$ws.Range("C277").Value = "La Araucanía"
$ws.Range("D277").Value = 44841
$ws.Range("E277").Value = 9
$ws.Range("F277").Value = 100112001
$ws.Range("G277").Value = "Berenjena"
$ws.Range("H277").Value = "Sin especificar"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 50
$ws.Range("K277").Value = 15000
$ws.Range("L277").Value = 15000
$ws.Range("M277").Value = 15000
$ws.Range("N277").Value = "$/caja 40 unidades"
$ws.Range("O277").Value = "Región de Arica y Parinacota"
$ws.Range("P277").Value = 375
$ws.Range("Q277").Value = 40
$ws.Range("R277").Value = "Hortaliza"
